$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values with newly calculated results
$ws.Range("C2").Value = 12919.98153166996

$ws.Range("B3").Value = 542.5
$ws.Range("C3").Value = 12847.20382482663

$ws.Range("B4").Value = 575
$ws.Range("C4").Value = 12814.35181620871

# Remove the last data row (row 5) entirely
$ws.Rows("5:5").Delete()
